$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix the typo "εκτλεί" -> "εκτελεί"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("εκτλεί", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "εκτελεί", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Remove one of the two consecutive, identical empty paragraphs that sit
#    right after "& ΜΗΧΑΝΙΚΩΝ ΥΠΟΛΟΓΙΣΤΩΝ" (keeps the one that sits right
#    before the paragraph carrying themeTint/themeShade on its color).
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Start = 0
$anchor.End = $d.Content.End
$anchor.Find.Execute("& ΜΗΧΑΝΙΚΩΝ ΥΠΟΛΟΓΙΣΤΩΝ", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null

$paraCount = $d.Paragraphs.Count
$hostIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if (($cand.Range.Start -le $anchor.Start) -and ($cand.Range.End -ge $anchor.End)) {
        $hostIndex = $i
    }
}
$dupEmptyPara = $d.Paragraphs.Item($hostIndex + 1)
$dupEmptyPara.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3) "Διάγραμμα περίπτωσης χρήσης (Βασική Ροή και ολοκλήρωση ανάλυσης )."
#    -> "Διάγραμμα ακολουθίας (Βασική Ροή και ολοκλήρωση ανάλυσης )."
#    (typed as "ακολουθία" then "ς" was appended, so it lands in 3 runs)
# ---------------------------------------------------------------------------
$sentence1 = $d.Content
$sentence1.Start = 0
$sentence1.End = $d.Content.End
$sentence1.Find.Execute("περίπτωσης χρήσης (Βασική", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$sentence1Start = $sentence1.Start

$w1 = $d.Content
$w1.Start = $sentence1Start
$w1.End = $sentence1Start + 17
$w1.Find.Execute("περίπτωσης χρήσης", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "ακολουθία", 2) | Out-Null

$w2 = $d.Range($w1.End, $w1.End)
$w2.InsertAfter("ς")

# ---------------------------------------------------------------------------
# 4) "Διάγραμμα περίπτωσης χρήσης (Εναλλακτική ροή Β και εισαγωγή dataset)."
#    -> "Διάγραμμα ακολουθίας (Εναλλακτική ροή Β και εισαγωγή dataset)."
# ---------------------------------------------------------------------------
$sentence2 = $d.Content
$sentence2.Start = 0
$sentence2.End = $d.Content.End
$sentence2.Find.Execute("περίπτωσης χρήσης (Εναλλακτική", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$sentence2Start = $sentence2.Start

$x1 = $d.Content
$x1.Start = $sentence2Start
$x1.End = $sentence2Start + 18
$x1.Find.Execute("περίπτωσης χρήσης", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "ακολουθίας", 2) | Out-Null

Write-Output "done"
